# Weekly price update: a new daily-price record is inserted as row 719,
# pushing every existing record from the old row 719 onward down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 719 (shifts rows 719:761 down to 720:762).
$ws.Rows.Item(719).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A719").Value = 4
$ws.Range("B719").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C719").Value = "Los Lagos"
$ws.Range("D719").Value = 44931
$ws.Range("E719").Value = 10
$ws.Range("F719").Value = 100112033
$ws.Range("G719").Value = "Lechuga"
$ws.Range("H719").Value = "Escarola"
$ws.Range("I719").Value = "Primera"
$ws.Range("J719").Value = 300
$ws.Range("K719").Value = 15000
$ws.Range("L719").Value = 16000
$ws.Range("M719").Value = 15500
$ws.Range("N719").Value = "`$/caja 15 unidades"
$ws.Range("O719").Value = "Región de Coquimbo"
$ws.Range("P719").Value = 1033
$ws.Range("Q719").Value = 15
$ws.Range("R719").Value = "Hortaliza"
